# Updates the cryptos list (Price / Volume(1h) columns, and the
# ordi/FraxShare row swap) to match the latest GitHub Actions refresh.
# Numeric-looking "Price" strings are assigned with a leading apostrophe
# so Excel keeps them as text (matching the original inlineStr cells)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.044.75'
$ws.Range("E2").Value = '  -1.70%  '
$ws.Range("D3").Value = '2.302.61'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''317.07'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = '''104.85'
$ws.Range("E6").Value = '  -2.62%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").Value = '''39.75'
$ws.Range("E10").Value = '  -4.44%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Value = '''0.107'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = '''0.979'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("E15").Value = '  -3.43%  '
$ws.Range("D16").Value = '2.650.74'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").Value = '2.300.00'
$ws.Range("E17").Value = '  -2.50%  '
$ws.Range("D18").Value = '42.009.18'
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("D19").Value = '''7.78'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").Value = '''286.73'
$ws.Range("E21").Value = '  +11.44%  '
$ws.Range("D22").Value = '''73.68'
$ws.Range("E22").Value = '  -3.93%  '
$ws.Range("D23").Value = '''3.57'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("D25").Value = '''9.99'
$ws.Range("E25").Value = '  +6.08%  '
$ws.Range("D26").Value = '''1.01'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("E27").Value = '  -3.93%  '
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").Value = '''2.22'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").Value = '''164.97'
$ws.Range("E30").Value = '  -5.66%  '
$ws.Range("D31").Value = '''35.49'
$ws.Range("E31").Value = '  -2.96%  '
$ws.Range("D32").Value = '''0.0884'
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("E34").Value = '  -3.12%  '
$ws.Range("D35").Value = '''0.132'
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").Value = '''0.116'
$ws.Range("E36").Value = '  -8.93%  '
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  +9.97%  '
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("E40").Value = '  -4.87%  '
$ws.Range("D41").Value = '''103.26'
$ws.Range("E41").Value = '  +21.36%  '
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("D43").Value = '''71.49'
$ws.Range("E44").Value = '  -5.27%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").Value = '''116.30'
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("D47").Value = '''12.15'
$ws.Range("E47").Value = '  +0.91%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '''9.18'
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = '''78.30'
$ws.Range("E49").Value = '  +4.08%  '
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("D51").Value = '''1.29'
$ws.Range("E51").Value = '  +1.86%  '
